$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the hyperlink object attached to E19 (the "USB Block" alibaba link)
# before we blank the row, so it doesn't linger as an orphan hyperlink.
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$E$19') {
        $hl.Delete()
    }
}

# --- Row 6 (J1 header pin): quantity corrected from 4 to 1
$ws.Range("F6").Value = 1

# --- Row 17 (Real-Time Clock & Battery):
#     add the new reference link and update the sourcing note for Production Run 3
$ws.Range("E17").Value = "https://www.amazon.com/HiLetgo-DS3231-Precision-Arduino-Raspberry/dp/B01N1LZSK3"
$ws.Range("G17").Value = "Can use same as last production, requirements: DS3231 RTC, 4 Female Pin With Battery. Please assemble onto PCB (unlike last production)"

# --- Row 19 ("USB Block" alternate-part line item): no longer needed, clear it out
$ws.Range("A19:G19").ClearContents()

# --- Move the active selection to C21 (matches author's final cursor position)
$ws.Range("C21").Select()
